$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item("Tableau2")

$formula = "=Tableau2[[#This Row],[Fin]]-Tableau2[[#This Row],[Début]]"

# Complete row 4 (Fin was missing) and extend the table with 7 new rows
$ws.Range("C4").Value = 0.5

for ($i = 0; $i -lt 7; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# Row 5
$ws.Range("A5").Value = 44683
$ws.Range("B5").Value = 0.5625
$ws.Range("C5").Value = 0.625

# Row 6
$ws.Range("A6").Value = 44683
$ws.Range("B6").Value = 0.63888888888888895
$ws.Range("C6").Value = 0.67013888888888884

# Row 7
$ws.Range("A7").Value = 44684
$ws.Range("B7").Value = 0.33333333333333331

# Recompute the "Total" calculated column for every row in the table range
for ($r = 4; $r -le 11; $r++) {
    $ws.Range("D$r").Formula = $formula
}

# Shared-text columns, written in the same order they first appear
$ws.Range("E4").Value = "création de la documentation de projet"
$ws.Range("F4").Value = "intro, objectif ajouter dans la doc"
$ws.Range("E5").Value = "planification initiale"
$ws.Range("E6").Value = "préparatin des logiciels a utiliser "
$ws.Range("F6").Value = "phpStorm, MySQL workbench, HeidiSQL"
$ws.Range("E7").Value = "création du projet sur phpStorm"
$ws.Range("F7").Value = "index & structure MVC"

$ws.Range("F9").Select() | Out-Null
